# Add a bottom-right "PAGE" field to the document's default footer.
$d = $word.ActiveDocument

# The primary (default) footer of the first/only section.
$sec = $d.Sections(1)
$ftr = $sec.Footers(1)
$rng = $ftr.Range

# Build the footer paragraph exactly as Word would emit it: right-aligned,
# "Footer" style, Calibri 9pt (sz=18 half-points) runs wrapping a classic
# begin/instrText/end PAGE field.
$footerXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="Footer"/><w:jc w:val="right"/></w:pPr>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="18"/></w:rPr>' +
  '<w:fldChar w:fldCharType="begin"/></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="18"/></w:rPr>' +
  '<w:instrText xml:space="preserve"> PAGE </w:instrText></w:r>' +
  '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:sz w:val="18"/></w:rPr>' +
  '<w:fldChar w:fldCharType="end"/></w:r>' +
  '</w:p>'

[void]$rng.InsertXML($footerXml)

Write-Output "Footer page number inserted."
